# Adds three new rows to the statistics table, right after the "GAD-7"
# row and before the trailing footnote row ("n (%); Mean (SD) (Min-Max)"):
#   - "Igi (VD)"        (bold section header, empty value cell)
#   - "    Insone"       / "30 (48%)"
#   - "    Não Insone"   / "33 (52%)"
#
# This mirrors the existing "Sexo" / "Feminino" / "Masculino" group that
# is already present in the table.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$W_NS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$borderColorInt = [System.Convert]::ToInt32("D3D3D3", 16)

function Set-CellBorders($cell) {
    foreach ($idx in @(-1, -2, -3, -4)) {
        $cell.Borders.Item($idx).LineStyle = 1
        $cell.Borders.Item($idx).Color = $borderColorInt
    }
    $cell.Borders.DistanceFromTop = 0
    $cell.Borders.DistanceFromBottom = 0
    $cell.Borders.DistanceFromLeft = 0
    $cell.Borders.DistanceFromRight = 0
}

function Set-CellParagraph($cell, [string]$xmlFragment) {
    # InsertXML on a freshly-split / still-empty cell behaves oddly when the
    # cell's Range is zero-length (collapsed): so first stuff a one-char
    # placeholder in to give the range real width, then shrink the
    # replacement range by one unit to avoid touching the end-of-cell mark
    # (doing InsertXML across that boundary duplicates the paragraph).
    $cell.Range.Text = "X"
    $target = $d.Range($cell.Range.Start, $cell.Range.End - 1)
    $target.InsertXML($xmlFragment)
}

function Add-StatRow([string]$label, [string]$value, [bool]$labelBold) {
    $footnoteRow = $t.Rows.Item($t.Rows.Count)
    $newRow = $t.Rows.Add($footnoteRow)
    $newRow.AllowBreakAcrossPages = $false

    $cell1 = $newRow.Cells.Item(1)
    $cell1.Split(1, 2)
    $cell1 = $newRow.Cells.Item(1)
    $cell2 = $newRow.Cells.Item(2)

    Set-CellBorders $cell1
    Set-CellBorders $cell2

    if ($labelBold) {
        $boldTag = '<w:b w:val="true"/>'
    } else {
        $boldTag = ''
    }

    $labelXml = "<w:p $W_NS><w:pPr><w:spacing w:before=`"0`" w:after=`"60`"/><w:keepNext/><w:jc w:val=`"start`"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Calibri`" w:hAnsi=`"Calibri`"/><w:sz w:val=`"20`"/>$boldTag</w:rPr><w:t xml:space=`"default`">$label</w:t></w:r></w:p>"
    $valueXml = "<w:p $W_NS><w:pPr><w:spacing w:before=`"0`" w:after=`"60`"/><w:keepNext/><w:jc w:val=`"center`"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=`"Calibri`" w:hAnsi=`"Calibri`"/><w:sz w:val=`"20`"/></w:rPr><w:t xml:space=`"default`">$value</w:t></w:r></w:p>"

    Set-CellParagraph $cell1 $labelXml
    Set-CellParagraph $cell2 $valueXml

    return $newRow
}

# Rows are added in reverse order since each Add() inserts immediately
# before the footnote row; doing the last-desired row first keeps the
# final on-page order: Igi (VD) / Insone / Não Insone.
Add-StatRow "    Não Insone" "33 (52%)" $false | Out-Null
Add-StatRow "    Insone" "30 (48%)" $false | Out-Null
Add-StatRow "Igi (VD)" "" $true | Out-Null

Write-Host "Final row count:" $t.Rows.Count
